{"js": "// Update the introduction text in the template (About this document section).\n//\n// 1) \"...data fields, types and descriptions...\" / \"...description in each\n//    field outlines the context for the data field, the frequency and\n//    milestone...\" -> \"...data fields, types, and descriptions...\" /\n//    \"...description for each field outlines the context, collection\n//    frequency and milestone...\"\n// 2) \"...noted at the end of the document...\" / \"The image below outlines\n//    the key milestones...\" -> \"...noted in the grey tables for each\n//    record...\" / \"The images below outline the key milestones...\"\n\nconst body = context.document.body;\n\nconst before1 =\n  \"The document outlines the data fields, types and descriptions for all required data collection during the lifetime of the Refugee Transition Outcome Fund. The description in each field outlines the context for the data field, the frequency and milestone at which it should be collected. \";\nconst after1 =\n  \"The document outlines the data fields, types, and descriptions for all required data collection during the lifetime of the Refugee Transition Outcome Fund. The description for each field outlines the context, collection frequency and milestone at which it should be collected. \";\n\nconst before2 =\n  \"The validation rules that apply for each field are noted at the end of the document. These will be applied on submission of data, if any validation rules are not met, a report will be produced with details of the fields and specified errors for the given fields.  The image below outlines the key milestones and data collection points during the program.\";\nconst after2 =\n  \"The validation rules that apply for each field are noted in the grey tables for each record. These will be applied on submission of data, if any validation rules are not met, a report will be produced with details of the fields and specified errors for the given fields.  The images below outline the key milestones and data collection points during the program.\";\n\nconst results1 = body.search(before1, { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\nfor (const range of results1.items) {\n  // clear() first so the replacement collapses into a single run and the\n  // stale grammar-check markers that used to split it are dropped too.\n  range.clear();\n  range.insertText(after1, \"Replace\");\n}\nawait context.sync();\n\nconst results2 = body.search(before2, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\nfor (const range of results2.items) {\n  range.clear();\n  range.insertText(after2, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the introduction text in the template (About this document section).\n#\n# 1) \"...data fields, types and descriptions...\" / \"...description in each\n#    field outlines the context for the data field, the frequency and\n#    milestone...\" -> \"...data fields, types, and descriptions...\" /\n#    \"...description for each field outlines the context, collection\n#    frequency and milestone...\"\n# 2) \"...noted at the end of the document...\" / \"The image below outlines\n#    the key milestones...\" -> \"...noted in the grey tables for each\n#    record...\" / \"The images below outline the key milestones...\"\n\n$d = $word.ActiveDocument\n\n$before1 = \"The document outlines the data fields, types and descriptions for all required data collection during the lifetime of the Refugee Transition Outcome Fund. The description in each field outlines the context for the data field, the frequency and milestone at which it should be collected. \"\n$after1  = \"The document outlines the data fields, types, and descriptions for all required data collection during the lifetime of the Refugee Transition Outcome Fund. The description for each field outlines the context, collection frequency and milestone at which it should be collected. \"\n\n$before2 = \"The validation rules that apply for each field are noted at the end of the document. These will be applied on submission of data, if any validation rules are not met, a report will be produced with details of the fields and specified errors for the given fields.  The image below outlines the key milestones and data collection points during the program.\"\n$after2  = \"The validation rules that apply for each field are noted in the grey tables for each record. These will be applied on submission of data, if any validation rules are not met, a report will be produced with details of the fields and specified errors for the given fields.  The images below outline the key milestones and data collection points during the program.\"\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute($before1, $false, $false, $false, $false, $false, $true, 1, $false, $after1, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute($before2, $false, $false, $false, $false, $false, $true, 1, $false, $after2, 2) | Out-Null\n"}
